$d = $word.ActiveDocument

# 1) Programa em português (paragraph split into runs separated by manual line breaks)
$d.Content.Find.Execute("Programa em português1.Conceitos básicos de Química (2 horas)a.Estrutura Atômicab.Tabela Periódicac.Ligações Químicas2.Os estados físicos da matéria e suas propriedades peculiares (6 horas)a.O estado gasoso – pressão, relações PVT, gases ideais e reaisb.O estado líquido – soluções, forças intermoleculares, viscosidade, tensão superficial, pressão de vapor, mudanças de fasec.O estado sólido – classificação dos sólidos (moleculares, reticulares, metálicos e iônicos) 3.Reações químicas (8 horas)a.Tipos de reações (dupla-troca, oxirredução)b.Estequiometria em reações químicas (reagentes limitantes, pureza e rendimento)c.Energia e reações químicasd.Equilíbrio químico – soluções tampãoe.Fundamentos de corrosão4.Noções de química orgânica (6 horas)a.Hidrocarbonetos e suas principais propriedadesb.Combustíveis e combustãoc.Polímeros5.Tecnologia Química aplicada (8 horas)a.Papel e celuloseb.Açúcar e álcoolc.Sabões e detergentesd.Petróleo e gáse.Gases industriais f.Produção de vidros e cimento", $true, $false, $false, $false, $false, $true, 1, $false, "Programa em português^l1.Conceitos básicos de Química (2 horas)^la.Estrutura Atômica^lb.Tabela Periódica^lc.Ligações Químicas^l2.Os estados físicos da matéria e suas propriedades peculiares (6 horas)^la.O estado gasoso – pressão, relações PVT, gases ideais e reais^lb.O estado líquido – soluções, forças intermoleculares, viscosidade, tensão superficial, pressão de vapor, mudanças de fase^lc.O estado sólido – classificação dos sólidos (moleculares, reticulares, metálicos e iônicos) ^l3.Reações químicas (8 horas)^la.Tipos de reações (dupla-troca, oxirredução)^lb.Estequiometria em reações químicas (reagentes limitantes, pureza e rendimento)^lc.Energia e reações químicas^ld.Equilíbrio químico – soluções tampão^le.Fundamentos de corrosão^l4.Noções de química orgânica (6 horas)^la.Hidrocarbonetos e suas principais propriedades^lb.Combustíveis e combustão^lc.Polímeros^l5.Tecnologia Química aplicada (8 horas)^la.Papel e celulose^lb.Açúcar e álcool^lc.Sabões e detergentes^ld.Petróleo e gás^le.Gases industriais ^lf.Produção de vidros e cimento", 2)

# 2) Programa em inglês (italic run)
$d.Content.Find.Execute("1.Basic Concepts of Chemistrya.Atomic Structureb.Periodic tablec.Chemical bonds2.The physical states of matter and their peculiar propertiesa.The gaseous state – pressure, PVT relations, ideal and real gasesb.The liquid state - solutions, intermolecular forces, viscosity, surface tension, vapor pressure, phase changesc.The solid state - classification of solids (molecular, reticular, metallic and ionic)3.Chemical reactionsa.Types of reactions (double-exchange, oxy-reduction)b.Stoichiometry in chemical reactions (limiting reagents, purity and yield) c.Energy and chemical reactionsd.Corrosion Fundamentals4.Notions of organic chemistrya.Hydrocarbons and their main propertiesb.Fuel and combustionc.Polymers5.Applied Chemistry Technologya.Paper and Celluloseb.Sugar and alcoholc.Soaps and detergentsd.Oil and gase.Industrial gasesf.Glass and cement production", $true, $false, $false, $false, $false, $true, 1, $false, "1.Basic Concepts of Chemistry^la.Atomic Structure^lb.Periodic table^lc.Chemical bonds^l2.The physical states of matter and their peculiar properties^la.The gaseous state – pressure, PVT relations, ideal and real gases^lb.The liquid state - solutions, intermolecular forces, viscosity, surface tension, vapor pressure, phase changes^lc.The solid state - classification of solids (molecular, reticular, metallic and ionic)^l3.Chemical reactions^la.Types of reactions (double-exchange, oxy-reduction)^lb.Stoichiometry in chemical reactions (limiting reagents, purity and yield) ^lc.Energy and chemical reactions^ld.Corrosion Fundamentals^l4.Notions of organic chemistry^la.Hydrocarbons and their main properties^lb.Fuel and combustion^lc.Polymers^l5.Applied Chemistry Technology^la.Paper and Cellulose^lb.Sugar and alcohol^lc.Soaps and detergents^ld.Oil and gas^le.Industrial gases^lf.Glass and cement production", 2)

# 3) Bibliografia
$d.Content.Find.Execute("BROWN, T.L. et al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007ATKINS, P. Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006KOTZ, J. C. et al. Química geral e reações químicas, 9ª Edição, São Paulo, Cengage Learning, 2015.TOLENTINO, N. M. C. Processos Químicos Industriais, 1ª Edição, São Paulo, Érica, 2015.", $true, $false, $false, $false, $false, $true, 1, $false, "BROWN, T.L. et al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007^lATKINS, P. Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006^lKOTZ, J. C. et al. Química geral e reações químicas, 9ª Edição, São Paulo, Cengage Learning, 2015.^lTOLENTINO, N. M. C. Processos Químicos Industriais, 1ª Edição, São Paulo, Érica, 2015.", 2)
